$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 11
$ws.Range("H11").Value = 36.3
$ws.Range("I11").Value = 36.3
$ws.Range("K11").Value = 36.3
$ws.Range("M11").Value = 103.7
# Row 19
$ws.Range("H19").Value = 802
$ws.Range("I19").Value = 750.5
$ws.Range("K19").Value = 750.5
$ws.Range("M19").Value = -575.5
# Row 28
$ws.Range("H28").Value = 357857.44
$ws.Range("I28").Value = 556157.25
$ws.Range("K28").Value = 556157.25
$ws.Range("M28").Value = -555672.25
# Row 33
$ws.Range("H33").Value = 91818730
$ws.Range("I33").Value = 1250353.6
$ws.Range("J33").Value = 333334400
$ws.Range("K33").Value = 1250353.6
$ws.Range("L33").Value = 333334400
$ws.Range("M33").Value = -1250124.6
$ws.Range("N33").Value = -333334858
# Row 43
$ws.Range("H43").Value = 9500
$ws.Range("I43").Value = 9500
$ws.Range("K43").Value = 9500
$ws.Range("M43").Value = -9431
# Row 69
$ws.Range("H69").Value = 9992.333000000001
$ws.Range("J69").Value = 14979
$ws.Range("L69").Value = 44937
$ws.Range("N69").Value = -46685
# Row 72
$ws.Range("H72").Value = 9992.333000000001
$ws.Range("J72").Value = 14979
$ws.Range("L72").Value = 134811
$ws.Range("N72").Value = -143547
# Row 86
$ws.Range("H86").Value = 2409.5
$ws.Range("I86").Value = 1500
$ws.Range("J86").Value = 2864.25
$ws.Range("K86").Value = 1500
$ws.Range("L86").Value = 2864.25
$ws.Range("M86").Value = -377
$ws.Range("N86").Value = -5110.25
# Row 89
$ws.Range("H89").Value = 2409.5
$ws.Range("I89").Value = 1500
$ws.Range("J89").Value = 2864.25
$ws.Range("K89").Value = 7500
$ws.Range("L89").Value = 14321.25
$ws.Range("M89").Value = -1884
$ws.Range("N89").Value = -25553.25
# Row 99
$ws.Range("H99").Value = 3706.1667
$ws.Range("I99").Value = 1161.3334
$ws.Range("J99").Value = 6251
$ws.Range("K99").Value = 3484.0002
$ws.Range("L99").Value = 18753
$ws.Range("M99").Value = -1986.0002
$ws.Range("N99").Value = -21749
# Row 100
$ws.Range("H100").Value = 2485.3635
$ws.Range("I100").Value = 2166.375
$ws.Range("K100").Value = 2166.375
$ws.Range("M100").Value = -1625.375
# Row 106
$ws.Range("H106").Value = 68065
$ws.Range("I106").Value = 81706.375
$ws.Range("J106").Value = 13499.5
$ws.Range("K106").Value = 81706.375
$ws.Range("L106").Value = 13499.5
$ws.Range("M106").Value = -81075.375
$ws.Range("N106").Value = -14761.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3301.3394
$ws.Range("I32").Value = 1650.3478
$ws.Range("K32").Value = 1650.3478
$ws.Range("M32").Value = -1363.3478
# Row 61
$ws.Range("H61").Value = 5234.2666
$ws.Range("I61").Value = 5250.5
$ws.Range("J61").Value = 5228.364
$ws.Range("K61").Value = 5250.5
$ws.Range("L61").Value = 5228.364
$ws.Range("M61").Value = -5038.5
$ws.Range("N61").Value = -5652.364
# Row 136
$ws.Range("H136").Value = 5234.2666
$ws.Range("I136").Value = 5250.5
$ws.Range("J136").Value = 5228.364
$ws.Range("K136").Value = 15751.5
$ws.Range("L136").Value = 15685.092
$ws.Range("M136").Value = -13201.5
$ws.Range("N136").Value = -20785.092

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 1800.2667
$ws.Range("I86").Value = 2051
$ws.Range("K86").Value = 2051
$ws.Range("M86").Value = -928
# Row 89
$ws.Range("H89").Value = 1800.2667
$ws.Range("I89").Value = 2051
$ws.Range("K89").Value = 10255
$ws.Range("M89").Value = -4639
# Row 134
$ws.Range("H134").Value = 3095.5312
$ws.Range("I134").Value = 3303.2273
$ws.Range("J134").Value = 2638.6
$ws.Range("K134").Value = 9909.6819
$ws.Range("L134").Value = 7915.799999999999
$ws.Range("M134").Value = -7374.6819
$ws.Range("N134").Value = -12985.8

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 107
$ws.Range("H107").Value = 1350.9
$ws.Range("I107").Value = 1003.8
$ws.Range("K107").Value = 1003.8
$ws.Range("M107").Value = 916.2
# Row 122
$ws.Range("H122").Value = 2159.5625
$ws.Range("I122").Value = 2131.6
$ws.Range("J122").Value = 2206.1667
$ws.Range("K122").Value = 6394.799999999999
$ws.Range("L122").Value = 6618.500100000001
$ws.Range("M122").Value = -3944.799999999999
$ws.Range("N122").Value = -11518.5001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 483.8
$ws.Range("I5").Value = 462
$ws.Range("J5").Value = 498.33334
$ws.Range("K5").Value = 1386
$ws.Range("L5").Value = 1495.00002
$ws.Range("M5").Value = -1274
$ws.Range("N5").Value = -1719.00002
# Row 135
$ws.Range("H135").Value = 483.8
$ws.Range("I135").Value = 462
$ws.Range("J135").Value = 498.33334
$ws.Range("K135").Value = 4158
$ws.Range("L135").Value = 4485.00006
$ws.Range("M135").Value = -1623
$ws.Range("N135").Value = -9555.00006

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 26454.482
$ws.Range("J70").Value = 10164.25
$ws.Range("L70").Value = 10164.25
$ws.Range("N70").Value = -10704.25
# Row 73
$ws.Range("H73").Value = 26454.482
$ws.Range("J73").Value = 10164.25
$ws.Range("L73").Value = 10164.25
$ws.Range("N73").Value = -12036.25
# Row 80
$ws.Range("H80").Value = 11000
$ws.Range("I80").Value = 20000
$ws.Range("K80").Value = 20000
$ws.Range("M80").Value = -19002
# Row 83
$ws.Range("H83").Value = 11000
$ws.Range("I83").Value = 20000
$ws.Range("K83").Value = 100000
$ws.Range("M83").Value = -95008
# Row 135
$ws.Range("H135").Value = 69902.8
$ws.Range("J135").Value = 74878.5
$ws.Range("L135").Value = 74878.5
$ws.Range("N135").Value = -85018.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 55
$ws.Range("H55").Value = 786.4666999999999
$ws.Range("I55").Value = 311.33334
$ws.Range("J55").Value = 1499.1666
$ws.Range("K55").Value = 311.33334
$ws.Range("L55").Value = 1499.1666
$ws.Range("M55").Value = -138.33334
$ws.Range("N55").Value = -1845.1666
# Row 100
$ws.Range("H100").Value = 3550.0908
$ws.Range("I100").Value = 3624.8333
$ws.Range("J100").Value = 3460.4
$ws.Range("K100").Value = 3624.8333
$ws.Range("L100").Value = 3460.4
$ws.Range("M100").Value = -3083.8333
$ws.Range("N100").Value = -4542.4
# Row 132
$ws.Range("H132").Value = 3435.568
$ws.Range("I132").Value = 3355.4055
$ws.Range("K132").Value = 10066.2165
$ws.Range("M132").Value = -7536.216499999999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 1269
$ws.Range("I132").Value = 1295.174
$ws.Range("J132").Value = 1068.3334
$ws.Range("K132").Value = 3885.522
$ws.Range("L132").Value = 3205.0002
$ws.Range("M132").Value = -1355.522
$ws.Range("N132").Value = -8265.0002
